# Ran Prod Verification Script
# Updates the "DateProd" (column B) timestamps written by the latest
# verification run across the three sheets: CC-Payments-Sale,
# CC-Payments-Auth and ACH-Payments-Debit.

$wb = $excel.ActiveWorkbook

$wsSale = $wb.Worksheets.Item("CC-Payments-Sale")
$wsSale.Range("B2").Value = "Thu Aug 28 07:50:38 IST 2025"

$wsAuth = $wb.Worksheets.Item("CC-Payments-Auth")
$wsAuth.Range("B2").Value = "Thu Aug 28 07:38:05 IST 2025"
$wsAuth.Range("B3").Value = "Thu Aug 28 07:38:48 IST 2025"
$wsAuth.Range("B4").Value = "Thu Aug 28 07:39:33 IST 2025"
$wsAuth.Range("B5").Value = "Thu Aug 28 07:40:39 IST 2025"
$wsAuth.Range("B6").Value = "Thu Aug 28 07:41:53 IST 2025"
$wsAuth.Range("B7").Value = "Thu Aug 28 07:42:54 IST 2025"

$wsDebit = $wb.Worksheets.Item("ACH-Payments-Debit")
$wsDebit.Range("B2").Value = "Thu Aug 28 07:43:40 IST 2025"
$wsDebit.Range("B3").Value = "Thu Aug 28 07:44:24 IST 2025"
$wsDebit.Range("B4").Value = "Thu Aug 28 07:45:10 IST 2025"
$wsDebit.Range("B5").Value = "Thu Aug 28 07:45:55 IST 2025"
$wsDebit.Range("B6").Value = "Thu Aug 28 07:46:47 IST 2025"
$wsDebit.Range("B7").Value = "Thu Aug 28 07:47:34 IST 2025"
$wsDebit.Range("B8").Value = "Thu Aug 28 07:48:18 IST 2025"
$wsDebit.Range("B9").Value = "Thu Aug 28 07:49:02 IST 2025"
$wsDebit.Range("B10").Value = "Thu Aug 28 07:49:45 IST 2025"
